# Update the "F" column (attendee/views count) figures on the "展览"
# (Worksheets index 1) and "全部类型" (Worksheets index 4) sheets.
# Both sheets list the same events, so the same value lands on each,
# just at slightly different row numbers because "全部类型" has a few
# extra rows interleaved.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# Row => new F value, for the "展览" sheet
$exhibitionUpdates = @{
    3  = 598
    4  = 301
    7  = 781
    9  = 202
    11 = 481
    12 = 1433
    13 = 136
    17 = 106
    18 = 686
    20 = 45
    21 = 277
    23 = 6175
    26 = 130
    28 = 14933
    29 = 1481
    31 = 117
    33 = 10885
    34 = 695
    35 = 4254
    36 = 194
}

# Row => new F value, for the "全部类型" sheet
$allTypesUpdates = @{
    3  = 598
    4  = 301
    7  = 781
    9  = 202
    11 = 481
    12 = 1433
    13 = 136
    18 = 106
    19 = 686
    22 = 45
    23 = 277
    26 = 6175
    29 = 130
    31 = 14933
    32 = 1481
    34 = 117
    36 = 10885
    37 = 695
    38 = 4254
    39 = 194
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
